$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 corresponds to the "GS + SB" strategy.
# D4 (total hits) changed from 442 to 489, which changes the derived
# precision (I4/L4) and F1-score (K4/N4) metrics accordingly.
$ws.Range("D4").Value = 489
$ws.Range("I4").Value = 0.02249488752556237
$ws.Range("K4").Value = 0.0437375745526839
$ws.Range("L4").Value = 0.02249488752556237
$ws.Range("N4").Value = 0.0437375745526839
